$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (19) down into the three
# new rows so the new cells inherit the same style (s="1", wrap text, etc.)
$ws.Range("A19:B19").Copy()
$ws.Range("A20:B22").PasteSpecial(-4122)

# Row 20 - 13 feb 11am
$ws.Range("A20").Value = "13 feb 11am"
$ws.Range("B20").Value = "Marginal model bayesian heterogeneity working. Way too slow speed."

# Row 21 - 13 Feb noon
$ws.Range("A21").Value = "13 Feb noon"
$ws.Range("B21").Value = "Marginal model bayesian heterogeneity has less autocorrelation. Takes less number of iterations to give the same results. Also it seems it gives more accurate results with very less label switching (which can be attributed to quick convergence)"
$ws.Rows.Item(21).RowHeight = 60

# Row 22 - 13 Feb 4:30pm
$ws.Range("A22").Value = "13 Feb 4:30pm"
$ws.Range("B22").Value = "Marginal model with unequal probabilities is working. But don't put the constraints that means should be ordered or sorted. Otherwise results are a nightmare"
$ws.Rows.Item(22).RowHeight = 45

# Update selection / view to match the new last cell, same as Excel would
# leave it after the edits were made.
$null = $ws.Range("B22").Select()

Write-Host "done"
